# Generate Report for Handoff
# - Update status text from "Handed back: in sync with en-US" to "Ready for handoff"
# - Bump the "Latest Handoff Datetime" / "Latest HO Xliff Generate Date" timestamps
# - Narrow the (now shorter) status column widths that used to fit the old text

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamps bumped forward (new report generation run) ---
$wsOverview.Range("G2").Value = "2016-09-03 17:07:32"
$wsDeDe.Range("H2").Value = "2016-09-03 17:07:32"
$wsZhCn.Range("H2").Value = "2016-09-03 17:07:23"

# --- Column widths shrink now that the status column holds shorter text ---
$wsOverview.Range("E1:F1").ColumnWidth = 16.333333333333332
$wsZhCn.Range("C1").ColumnWidth = 16.333333333333332
$wsDeDe.Range("C1").ColumnWidth = 16.333333333333332
